$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '91.904.30'
$ws.Range('E2').Value = '  +2.64%  '
$ws.Range('D3').Value = '3.165.51'
$ws.Range('E3').Value = '  +2.99%  '
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('D5').Value = '240.07'
$ws.Range('E5').Value = '  +2.17%  '
$ws.Range('D6').Value = '620.82'
$ws.Range('E6').Value = '  +0.52%  '
$ws.Range('D7').Value = '1.12'
$ws.Range('E7').Value = '  +6.23%  '
$ws.Range('D8').Value = '0.377'
$ws.Range('E8').Value = '  +4.51%  '
$ws.Range('D9').Value = '0.998'
$ws.Range('E9').Value = '  -0.18%  '
$ws.Range('D10').Value = '3.161.60'
$ws.Range('E10').Value = '  +2.91%  '
$ws.Range('D11').Value = '0.747'
$ws.Range('E11').Value = '  +4.80%  '
$ws.Range('E12').Value = '  +2.99%  '
$ws.Range('D13').Value = '0.0000248'
$ws.Range('E13').Value = '  -0.14%  '
$ws.Range('D14').Value = '35.57'
$ws.Range('E14').Value = '  +1.59%  '
$ws.Range('D15').Value = '5.57'
$ws.Range('E15').Value = '  +4.02%  '
$ws.Range('D16').Value = '91.748.07'
$ws.Range('E16').Value = '  +2.45%  '
$ws.Range('D17').Value = '3.745.91'
$ws.Range('E17').Value = '  +2.87%  '
$ws.Range('D18').Value = '3.145.85'
$ws.Range('E18').Value = '  +1.52%  '
$ws.Range('D19').Value = '3.77'
$ws.Range('E19').Value = '  -0.78%  '
$ws.Range('D20').Value = '15.22'
$ws.Range('E20').Value = '  +10.95%  '
$ws.Range('D21').Value = '5.94'
$ws.Range('E21').Value = '  +10.18%  '
$ws.Range('B22').Value = 'BitcoinCash'
$ws.Range('C22').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D22').Value = '458.87'
$ws.Range('E22').Value = '  +6.18%  '
$ws.Range('B23').Value = 'PEPE'
$ws.Range('C23').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D23').Value = '0.0000203'
$ws.Range('E23').Value = '  -3.83%  '
$ws.Range('D24').Value = '9.23'
$ws.Range('E24').Value = '  +5.61%  '
$ws.Range('D25').Value = '6.06'
$ws.Range('E25').Value = '  +6.04%  '
$ws.Range('B26').Value = 'Litecoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D26').Value = '88.46'
$ws.Range('E26').Value = '  +1.96%  '
$ws.Range('B27').Value = 'Aptos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D27').Value = '12.08'
$ws.Range('E27').Value = '  +3.19%  '
$ws.Range('D28').Value = '3.319.51'
$ws.Range('E28').Value = '  +2.04%  '
$ws.Range('E29').Value = '  +0.21%  '
$ws.Range('B30').Value = 'Hedera'
$ws.Range('C30').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D30').Value = '0.127'
$ws.Range('E30').Value = '  +40.50%  '
$ws.Range('B31').Value = 'Stellar'
$ws.Range('C31').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D31').Value = '0.232'
$ws.Range('E31').Value = '  +19.42%  '
$ws.Range('B32').Value = 'Cronos'
$ws.Range('C32').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D32').Value = '0.172'
$ws.Range('E32').Value = '  +10.85%  '
$ws.Range('D33').Value = '9.40'
$ws.Range('E33').Value = '  +4.05%  '
$ws.Range('B34').Value = 'Kaspa'
$ws.Range('C34').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D34').Value = '0.171'
$ws.Range('E34').Value = '  +13.18%  '
$ws.Range('B35').Value = 'Binance-PegBSC-USD'
$ws.Range('C35').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D35').Value = '1.00'
$ws.Range('E35').Value = '  -0.03%  '
$ws.Range('D36').Value = '7.68'
$ws.Range('E36').Value = '  +8.23%  '
$ws.Range('D37').Value = '26.58'
$ws.Range('E37').Value = '  +4.35%  '
$ws.Range('D38').Value = '513.06'
$ws.Range('E38').Value = '  +3.70%  '
$ws.Range('D39').Value = '1.37'
$ws.Range('E39').Value = '  +9.06%  '
$ws.Range('D40').Value = '1.94'
$ws.Range('E40').Value = '  +3.44%  '
$ws.Range('D41').Value = '0.452'
$ws.Range('E41').Value = '  +14.08%  '
$ws.Range('D42').Value = '3.81'
$ws.Range('E42').Value = '  +6.42%  '
$ws.Range('D43').Value = '3.47'
$ws.Range('E43').Value = '  -5.67%  '
$ws.Range('D44').Value = '22.22'
$ws.Range('E44').Value = '  +0.59%  '
$ws.Range('D46').Value = '159.77'
$ws.Range('E46').Value = '  +4.91%  '
$ws.Range('B47').Value = 'ARBITRUM'
$ws.Range('C47').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D47').Value = '0.717'
$ws.Range('E47').Value = '  +6.47%  '
$ws.Range('B48').Value = 'Stacks'
$ws.Range('C48').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D48').Value = '1.95'
$ws.Range('E48').Value = '  +5.73%  '
$ws.Range('D49').Value = '1.38'
$ws.Range('E49').Value = '  +6.36%  '
$ws.Range('B50').Value = 'Filecoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D50').Value = '4.45'
$ws.Range('E50').Value = '  +3.09%  '
$ws.Range('B51').Value = 'OKB'
$ws.Range('C51').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D51').Value = '44.15'
$ws.Range('E51').Value = '  -0.40%  '
